# Added Gracefull exit handlers
# - Update the saved selection on Sheet1 (user had moved on to H11 before
#   switching tabs).
# - Add a new "Sheet2" summarising a load-test run (active threads,
#   response time, TPS) with a TPS formula and a note about the graceful
#   exit condition, and leave it as the active/selected sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Remember where the user last clicked on Sheet1 before moving to Sheet2.
$ws1.Range("H11").Select() | Out-Null

# New sheet, inserted right after Sheet1.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Labels (bold, matching the style already used for Sheet1's row/column
# headers) -- written in reading order so shared-string ids line up with
# the source workbook.
$ws2.Range("A2").Value = "Active threads"
$ws2.Range("A2").Font.Bold = $true
$ws2.Range("A3").Value = "Response Time"
$ws2.Range("A3").Font.Bold = $true
$ws2.Range("A4").Value = "TPS"
$ws2.Range("A4").Font.Bold = $true

# Values + the TPS formula.
$ws2.Range("B2").Value = 4
$ws2.Range("B3").Value = 1000
$ws2.Range("B4").Formula = "=(B2*1000)/B3"

# Graceful-exit note next to the thread count.
$ws2.Range("C2").Value = "Response time goes beyond 2000> or we encounter error"

# Column A sized to fit the labels.
$ws2.Columns.Item(1).AutoFit() | Out-Null

# Leave the cursor where the author left it, and make Sheet2 the active tab.
$ws2.Range("B3").Select() | Out-Null
